# Updated symbol list on Tue Jan  3 10:51:28 UTC 2023 with GitHub Actions
# Applies the cryptos.xlsx price/volume/coin-row updates described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.44%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.30%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.37%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.57%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.673"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.13%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.255"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.94%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8494"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.21%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8577"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.78%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.98%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07087"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.03%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03263"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "13.74%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.32%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001529"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.55%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005937"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.01%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006027"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.43%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.527"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.83%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.216"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.35%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.84%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.56%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1315"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.58%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.496"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.09%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1409"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.19%"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04124"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.37%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.90%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004150"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.79%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.80%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "4.75%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03753"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.01%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1071"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.07%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002469"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.61%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003540"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-37.19%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008924"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.08%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005483"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.23%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.01%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-20.23%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-10.84%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
